# GHC7-waypoints.xlsx — add a "connected-to" waypoint-graph column.
#
# Before: A=Location ID, B/C=Horizontal/Vertical Position (TL-in),
#         D/E=Horizontal/Vertical Position (C-in) [formulas], F=Bubble
#         Diameter (in) header/values, G=Meter-to-Inch ratio (G1) + the
#         Bubble Diameter value in G2.
#
# After: a brand-new column is inserted at F ("connected-to"), so the old
#        F/G columns shift right to G/H (formulas referencing $G$1 become
#        $H$1 automatically via Excel's reference-adjusting column insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column F (Bubble Diameter). This
# shifts the old F->G and G->H, and Excel auto-rewrites the $G$1 formula
# refs in D/E to $H$1.
$ws.Columns("F:F").Insert()

# Match the new column F's display width to its neighbour (column E) so it
# doesn't sit at the generic default width.
$ws.Columns("F:F").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Header
$ws.Cells.Item(1, 6).Value = "connected-to"

# Waypoint connection data (row -> "connected-to" value). Plain integers
# (no separator) land as numeric cells; multi-id lists land as text.
$connections = @{
    2  = 2
    3  = "1;3;7"
    4  = "2;4"
    5  = "3;5;9"
    6  = "4;10"
    7  = 7
    8  = "2;6;8"
    9  = "7;11"
    10 = "4;11"
    11 = "5;12"
    12 = "8;9;12;13"
    13 = "10;11"
    14 = "11;17"
    15 = 15
    16 = "14;16;19;20"
    17 = "15;17"
    18 = "13;16;18;21"
    19 = "17;22"
    20 = "15;23"
    21 = "15;23"
    22 = "17;25"
    23 = "18;26"
    24 = "19;20;24"
    25 = "23;27"
    26 = "21;26;27"
    27 = "22;25"
    28 = "24;25;28"
    29 = "27;29"
    30 = "28;30;31"
    31 = 29
    32 = "29;32"
    33 = 31
}

# Row 28 is written last (its shared string ends up appended at the end of
# the table rather than in row order) to match how the original was
# authored -- every other row's "connected-to" text was entered first and
# row 28's was filled in / corrected afterwards.
for ($row = 2; $row -le 33; $row++) {
    if ($row -ne 28) {
        $ws.Cells.Item($row, 6).Value = $connections[$row]
    }
}
$ws.Cells.Item(28, 6).Value = $connections[28]

# Restore view state as closely as possible: scroll so D3 is the
# top-left visible cell, and leave the selection on F33 (last edited cell).
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F33").Select() | Out-Null
